$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.215.96"
$ws.Range("E2").Value = "  -0.34%  "

$ws.Range("D3").Value = "3.204.14"
$ws.Range("E3").Value = "  +0.47%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "607.77"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.90%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "155.98"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.11%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("D8").Value = "3.206.51"
$ws.Range("E8").Value = "  +0.63%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.552"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.00%  "

$ws.Range("E10").Value = "  -0.23%  "

$ws.Range("E11").Value = "  -4.39%  "

$ws.Range("E12").Value = "  -2.87%  "

$ws.Range("E13").Value = "  +0.73%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "38.51"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.81%  "

$ws.Range("D15").Value = "3.731.78"
$ws.Range("E15").Value = "  +0.53%  "

$ws.Range("D16").Value = "66.376.01"
$ws.Range("E16").Value = "  -0.13%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.34"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.12%  "

$ws.Range("D18").Value = "3.207.98"
$ws.Range("E18").Value = "  +0.37%  "

$ws.Range("E19").Value = "  +1.37%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "508.27"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.17%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.32"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.73%  "

$ws.Range("E22").Value = "  -0.49%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.02"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.15%  "

$ws.Range("E24").Value = "  -2.03%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "85.18"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.24%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.11%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.08%  "

$ws.Range("E28").Value = "  -2.12%  "

$ws.Range("E29").Value = "  +0.48%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.127"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +40.99%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.94"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.76%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.00"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.74%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "28.28"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.36%  "

$ws.Range("E34").Value = "  +0.09%  "

$ws.Range("E35").Value = "  -5.16%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.49"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.63%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "502.30"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.68%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "55.45"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.80%  "

$ws.Range("D39").Value = "0.0₃0775"
$ws.Range("E39").Value = "  +15.83%  "

$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.131"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.97%  "

$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0422"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.04%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.06"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +6.13%  "

$ws.Range("E43").Value = "  -2.17%  "

$ws.Range("E44").Value = "  -1.51%  "

$ws.Range("E45").Value = "  -0.54%  "

$ws.Range("D46").Value = "2.905.89"
$ws.Range("E46").Value = "  -0.02%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "28.44"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.08%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.40"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.92%  "

$ws.Range("E49").Value = "  +0.00%  "

$ws.Range("E50").Value = "  -1.69%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "122.50"
$ws.Range("D51").Style = "Normal"
